$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29: single label cell carried over from the "ab hier mit cuml"-style note row
$ws.Range("A29").Value = "ohne cuml "

# Row 30: new regular "full random" training run (plain/default formatting)
$ws.Range("A30").Value = "regular"
$ws.Range("B30").Value = "full random"
$ws.Range("C30").Value = "sentences"
$ws.Range("D30").Value = 1000
$ws.Range("E30").Value = 200
$ws.Range("F30").Value = 5
$ws.Range("H30").Value = "3000s"
$ws.Range("I30").Value = "null"
$ws.Range("J30").Value = "yes"
$ws.Range("K30").Value = 130
$ws.Range("L30").Value = "medium"

# Row 31: new zeroshot-huang-combined run, formatted like the other highlighted rows (row 16/19/21/22)
$ws.Range("A31").Value = "zeroshot huang combined with own (w/o marketing)"
$ws.Range("D31").Value = 1000
$ws.Range("E31").Value = 200
$ws.Range("F31").Value = 5
$ws.Range("G31").Value = 0.8
$ws.Range("H31").Value = "3300s"
$ws.Range("I31").Value = "null"
$ws.Range("J31").Value = "yes"
$ws.Range("K31").Value = 145
$ws.Range("L31").Value = "bad"
$ws.Range("N31").Value = "zeroshot minsim variert stark mit topic size!!!"

# Copy the existing "white fill" row format (row 16, A:L) onto the new row 31 (A:L)
# so it reuses the same cell style instead of minting a new one.
$ws.Range("A16:L16").Copy()
$ws.Range("A31:L31").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# N31 also belongs to the same highlighted row; copy the same format from a
# cell that already carries it (e.g. B16) onto N31.
$ws.Range("B16").Copy()
$ws.Range("N31").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move the selection to K31 (matches the post-edit saved view) and drop the old
# scrolled topLeftCell state.
$ws.Range("A1").Select()
$ws.Range("K31").Select()
